$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("04-13-2022")
$ws1.Delete()

$keySheet = $wb.Worksheets.Item("Key")
$keySheet.Activate()
$keySheet.Range("H14").Select()

$dat = $wb.Worksheets.Item("Daily Attendance Template")
$dat.Range("I6").ClearContents()
$dat.Range("I7").ClearContents()
$dat.Range("I8").Formula = "=I6-I7"

$dat.Protect($null, $true, $true, $true)

$dat.Activate()
$dat.Range("R1").Select()
